$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay stored as text,
# matching the source data (all Price/Volume cells are inline strings).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.585.72"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.987.99"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.38"
$ws.Range("E5").Value = "  +2.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.02"
$ws.Range("E6").Value = "  +4.23%  "

$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.979.00"
$ws.Range("E9").Value = "  +1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").Value = "  +3.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.39"
$ws.Range("E11").Value = "  +12.08%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +4.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.87"
$ws.Range("E14").Value = "  +3.31%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.476.16"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.09"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.982.71"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.513.76"
$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "439.22"
$ws.Range("E20").Value = "  +5.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.69"
$ws.Range("E21").Value = "  +2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("E22").Value = "  +3.96%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.04"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.13"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  +10.42%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.56"
$ws.Range("E29").Value = "  +2.56%  "

$ws.Range("E30").Value = "  +3.19%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.26"
$ws.Range("E31").Value = "  +5.03%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.86"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +10.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0781"
$ws.Range("E34").Value = "  +11.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.93"
$ws.Range("E35").Value = "  +4.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.982"
$ws.Range("E36").Value = "  +3.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.08"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.69"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.69"
$ws.Range("E39").Value = "  -3.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.78"
$ws.Range("E40").Value = "  +2.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.14"
$ws.Range("E41").Value = "  +4.02%  "

$ws.Range("E42").Value = "  +1.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.737.85"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  +6.23%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.12"
$ws.Range("E47").Value = "  +21.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.21"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.41"
$ws.Range("E51").Value = "  +2.19%  "
